# Scheduled market-data refresh: update currentAveragePrice* / Leve profit
# columns (H:N) for a set of rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR
# leve-profit tables to reflect freshly pulled Universalis prices.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 776.3158
$ws.Range("I6").Value = 106.89655
$ws.Range("J6").Value = 2933.3333
$ws.Range("K6").Value = 320.68965
$ws.Range("L6").Value = 8799.999899999999
$ws.Range("M6").Value = -208.68965
$ws.Range("N6").Value = -9023.999899999999
$ws.Range("H129").Value = 1157.075
$ws.Range("I129").Value = 586.25
$ws.Range("J129").Value = 1299.7812
$ws.Range("K129").Value = 1758.75
$ws.Range("L129").Value = 3899.3436
$ws.Range("M129").Value = 3241.25
$ws.Range("N129").Value = -13899.3436
$ws.Range("H138").Value = 2283.3572
$ws.Range("I138").Value = 3296.5
$ws.Range("J138").Value = 2007.0454
$ws.Range("K138").Value = 9889.5
$ws.Range("L138").Value = 6021.1362
$ws.Range("M138").Value = -4749.5
$ws.Range("N138").Value = -16301.1362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2518.879
$ws.Range("I61").Value = 1789.9474
$ws.Range("J61").Value = 3508.1428
$ws.Range("K61").Value = 1789.9474
$ws.Range("L61").Value = 3508.1428
$ws.Range("M61").Value = -1577.9474
$ws.Range("N61").Value = -3932.1428
$ws.Range("H132").Value = 4490.593
$ws.Range("I132").Value = 4732
$ws.Range("J132").Value = 4230.615
$ws.Range("K132").Value = 14196
$ws.Range("L132").Value = 12691.845
$ws.Range("M132").Value = -11666
$ws.Range("N132").Value = -17751.845
$ws.Range("H136").Value = 2518.879
$ws.Range("I136").Value = 1789.9474
$ws.Range("J136").Value = 3508.1428
$ws.Range("K136").Value = 5369.8422
$ws.Range("L136").Value = 10524.4284
$ws.Range("M136").Value = -2819.8422
$ws.Range("N136").Value = -15624.4284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2880.762
$ws.Range("I20").Value = 2374.4
$ws.Range("J20").Value = 4146.6665
$ws.Range("K20").Value = 2374.4
$ws.Range("L20").Value = 4146.6665
$ws.Range("M20").Value = -2127.4
$ws.Range("N20").Value = -4640.6665
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H86").Value = 62502108
$ws.Range("I86").Value = 71430550
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 71430550
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -71429427
$ws.Range("N86").Value = -5246
$ws.Range("H89").Value = 62502108
$ws.Range("I89").Value = 71430550
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 357152750
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -357147134
$ws.Range("N89").Value = -26232
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H130").Value = 25000
$ws.Range("J130").Value = 25000
$ws.Range("L130").Value = 25000
$ws.Range("N130").Value = -35040
$ws.Range("H133").Value = 98780
$ws.Range("J133").Value = 98780
$ws.Range("L133").Value = 98780
$ws.Range("N133").Value = -108900

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 298.8889
$ws.Range("I22").Value = 223.33333
$ws.Range("J22").Value = 450
$ws.Range("K22").Value = 223.33333
$ws.Range("L22").Value = 450
$ws.Range("M22").Value = 126.66667
$ws.Range("N22").Value = -1150
$ws.Range("H86").Value = 2551.25
$ws.Range("I86").Value = 2756.6875
$ws.Range("J86").Value = 2277.3333
$ws.Range("K86").Value = 2756.6875
$ws.Range("L86").Value = 2277.3333
$ws.Range("M86").Value = -1633.6875
$ws.Range("N86").Value = -4523.3333
$ws.Range("H89").Value = 2551.25
$ws.Range("I89").Value = 2756.6875
$ws.Range("J89").Value = 2277.3333
$ws.Range("K89").Value = 13783.4375
$ws.Range("L89").Value = 11386.6665
$ws.Range("M89").Value = -8167.4375
$ws.Range("N89").Value = -22618.6665
$ws.Range("H134").Value = 1872.9286
$ws.Range("I134").Value = 1810.0834
$ws.Range("K134").Value = 5430.2502
$ws.Range("M134").Value = -2895.2502

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1810.619
$ws.Range("I5").Value = 739.8
$ws.Range("J5").Value = 2145.25
$ws.Range("K5").Value = 2219.4
$ws.Range("L5").Value = 6435.75
$ws.Range("M5").Value = -2107.4
$ws.Range("N5").Value = -6659.75
$ws.Range("H7").Value = 376.15384
$ws.Range("J7").Value = 571.6667
$ws.Range("L7").Value = 1715.0001
$ws.Range("N7").Value = -1939.0001
$ws.Range("H34").Value = 11628297
$ws.Range("J34").Value = 13889340
$ws.Range("L34").Value = 41668020
$ws.Range("N34").Value = -41668188
$ws.Range("H39").Value = 1216.4
$ws.Range("I39").Value = 333.2
$ws.Range("J39").Value = 2099.6
$ws.Range("K39").Value = 999.5999999999999
$ws.Range("L39").Value = 6298.799999999999
$ws.Range("M39").Value = -705.5999999999999
$ws.Range("N39").Value = -6886.799999999999
$ws.Range("H40").Value = 234.11765
$ws.Range("I40").Value = 196.33333
$ws.Range("J40").Value = 276.625
$ws.Range("K40").Value = 785.33332
$ws.Range("L40").Value = 1106.5
$ws.Range("M40").Value = -716.33332
$ws.Range("N40").Value = -1244.5
$ws.Range("H55").Value = 1850
$ws.Range("J55").Value = 2187.5
$ws.Range("L55").Value = 6562.5
$ws.Range("N55").Value = -6916.5
$ws.Range("H80").Value = 5180
$ws.Range("J80").Value = 5300.3335
$ws.Range("L80").Value = 15901.0005
$ws.Range("N80").Value = -17773.0005
$ws.Range("H83").Value = 5180
$ws.Range("J83").Value = 5300.3335
$ws.Range("L83").Value = 47703.0015
$ws.Range("N83").Value = -57063.0015
$ws.Range("H122").Value = 5958.4736
$ws.Range("I122").Value = 515.3570999999999
$ws.Range("J122").Value = 21199.2
$ws.Range("K122").Value = 4638.2139
$ws.Range("L122").Value = 190792.8
$ws.Range("M122").Value = -2188.2139
$ws.Range("N122").Value = -195692.8
$ws.Range("H131").Value = 1007.7368
$ws.Range("I131").Value = 330
$ws.Range("J131").Value = 1019.8393
$ws.Range("K131").Value = 990
$ws.Range("L131").Value = 3059.5179
$ws.Range("M131").Value = 4050
$ws.Range("N131").Value = -13139.5179
$ws.Range("H132").Value = 3699.373
$ws.Range("I132").Value = 2630.1538
$ws.Range("J132").Value = 4377.4146
$ws.Range("K132").Value = 23671.3842
$ws.Range("L132").Value = 39396.7314
$ws.Range("M132").Value = -21141.3842
$ws.Range("N132").Value = -44456.7314
$ws.Range("H135").Value = 1810.619
$ws.Range("I135").Value = 739.8
$ws.Range("J135").Value = 2145.25
$ws.Range("K135").Value = 6658.2
$ws.Range("L135").Value = 19307.25
$ws.Range("M135").Value = -4123.2
$ws.Range("N135").Value = -24377.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6725
$ws.Range("I70").Value = 6800
$ws.Range("J70").Value = 6600
$ws.Range("K70").Value = 6800
$ws.Range("L70").Value = 6600
$ws.Range("M70").Value = -6530
$ws.Range("N70").Value = -7140
$ws.Range("H73").Value = 6725
$ws.Range("I73").Value = 6800
$ws.Range("J73").Value = 6600
$ws.Range("K73").Value = 6800
$ws.Range("L73").Value = 6600
$ws.Range("M73").Value = -5864
$ws.Range("N73").Value = -8472
$ws.Range("H113").Value = 2971
$ws.Range("J113").Value = 2971
$ws.Range("L113").Value = 2971
$ws.Range("N113").Value = -7311
$ws.Range("H132").Value = 3899.818
$ws.Range("I132").Value = 4000
$ws.Range("J132").Value = 3842.5715
$ws.Range("K132").Value = 12000
$ws.Range("L132").Value = 11527.7145
$ws.Range("M132").Value = -9470
$ws.Range("N132").Value = -16587.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4719.154
$ws.Range("I61").Value = 4999.875
$ws.Range("J61").Value = 4270
$ws.Range("K61").Value = 4999.875
$ws.Range("L61").Value = 4270
$ws.Range("M61").Value = -4797.875
$ws.Range("N61").Value = -4674
$ws.Range("H113").Value = 4719.154
$ws.Range("I113").Value = 4999.875
$ws.Range("J113").Value = 4270
$ws.Range("K113").Value = 4999.875
$ws.Range("L113").Value = 4270
$ws.Range("M113").Value = -2829.875
$ws.Range("N113").Value = -8610
$ws.Range("H136").Value = 2864.4614
$ws.Range("I136").Value = 3549.125
$ws.Range("J136").Value = 1769
$ws.Range("K136").Value = 10647.375
$ws.Range("L136").Value = 5307
$ws.Range("M136").Value = -8097.375
$ws.Range("N136").Value = -10407

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1819.5
$ws.Range("I126").Value = 1580.6666
$ws.Range("K126").Value = 4741.9998
$ws.Range("M126").Value = -2271.9998
